# feat(Import): Cantines: supprime le type de production "central" et met a jour les champs lies
#
# The "canteens_good" fixture drops its "nombre_satellites" column (K) and the
# siret value in A2 becomes a genuine number instead of a text-formatted string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (canteens_good -> canteens_good-2)
$ws.Name = "canteens_good-2"

# Iterative-calculation tolerance: 0.0001 -> 0.001
$excel.Iteration = $true
$excel.MaxChange = 0.001
$excel.MaxIterations = 100

# Column K ("nombre_satellites") is removed entirely.
$ws.Columns.Item(11).Delete()

# A2 ("siret") was stored as Text (format "@"); it becomes a plain number.
$ws.Range("A2").NumberFormat = "General"
$ws.Range("A2").Value = 21340172201787

# Drop the bespoke per-cell formatting the old sheet carried (extra fonts /
# text-number-format / protection xfs) so every cell falls back to the
# sheet's default style.
$ws.Cells.ClearFormats()

# Cosmetic view state that moved alongside the edit.
$ws.Range("A1").Select()
$excel.ActiveWindow.Zoom = 65

Write-Output "done"
